$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.327.95'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.691.06'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = "'218.41"
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = "'0.5281"
$ws.Range("E6").Value = '  +4.11%  '
$ws.Range("D7").Value = "'1.007"
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = "'0.2711"
$ws.Range("E8").Value = '  +2.06%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").Value = "'22.12"
$ws.Range("E9").Value = '  +2.94%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = "'0.06430"
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("D11").Value = "'0.07482"
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").Value = '1.713.20'
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("D13").Value = "'4.574"
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").Value = "'0.5863"
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").Value = "'0.000008521"
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").Value = "'64.56"
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = '26.384.03'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D18").Value = "'4.949"
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = "'10.91"
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").Value = "'189.47"
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = "'6.229"
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = "'144.93"
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").Value = "'7.713"
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").Value = "'0.1236"
$ws.Range("E26").Value = '  +5.26%  '
$ws.Range("D27").Value = "'15.88"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = "'0.06681"
$ws.Range("E28").Value = '  +15.01%  '
$ws.Range("D29").Value = "'1.358"
$ws.Range("E29").Value = '  +6.03%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").Value = "'3.588"
$ws.Range("E31").Value = '  +2.17%  '
$ws.Range("D32").Value = "'3.578"
$ws.Range("E32").Value = '  +1.26%  '
$ws.Range("D33").Value = "'1.672"
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("E34").Value = '  +2.05%  '
$ws.Range("D35").Value = "'0.6242"
$ws.Range("E35").Value = '  +4.22%  '
$ws.Range("D36").Value = "'2.394"
$ws.Range("E37").Value = '  +2.24%  '
$ws.Range("D38").Value = "'6.370"
$ws.Range("E38").Value = '  +5.88%  '
$ws.Range("D39").Value = '1.118.19'
$ws.Range("E39").Value = '  +4.41%  '
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").Value = "'0.8906"
$ws.Range("E41").Value = '  +3.20%  '
$ws.Range("D42").Value = "'1.018"
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("D43").Value = "'100.92"
$ws.Range("E43").Value = '  +1.38%  '
$ws.Range("D44").Value = '1.839.04'
$ws.Range("D45").Value = "'0.00000000115"
$ws.Range("E45").Value = '  +3.68%  '
$ws.Range("D46").Value = "'57.03"
$ws.Range("E46").Value = '  +2.41%  '
$ws.Range("D47").Value = "'8.184"
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = "'0.05269"
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").Value = "'6.118"
$ws.Range("E50").Value = '  +4.39%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = "'0.4303"
$ws.Range("E51").Value = '  +0.26%  '
